$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 85, shifting existing rows 85-107 down to 86-108.
$ws.Rows("85:85").Insert()

# Populate the newly inserted row 85 with the new weekly record.
$ws.Range("A85").Value = 1
$ws.Range("B85").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C85").Value = "Arica y Parinacota"
$ws.Range("D85").Value = 44754
$ws.Range("E85").Value = 15
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100106
$ws.Range("H85").Value = "Oleaginosos"
$ws.Range("I85").Value = 100106002
$ws.Range("J85").Value = "Palta"
$ws.Range("K85").Value = "Hass"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 400
$ws.Range("N85").Value = 15000
$ws.Range("O85").Value = 16000
$ws.Range("P85").Value = 15500
$ws.Range("Q85").Value = "$/bandeja 10 kilos"
$ws.Range("R85").Value = "Perú"
$ws.Range("S85").Value = 1550
$ws.Range("T85").Value = 10
